$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-28 Wednesday" "2026-01-29 Thursday"

Replace-Text "88×63=" "47×56="
Replace-Text "31×82=" "65×89="
Replace-Text "92×73=" "41×13="
Replace-Text "33×61=" "20×48="
Replace-Text "65×78=" "57×41="

Replace-Text "46×11=" "88×85="
Replace-Text "17×19=" "25×13="
Replace-Text "53×99=" "99×53="
Replace-Text "36×12=" "59×79="
Replace-Text "42×80=" "95×68="

Replace-Text "34×48=" "64×13="
Replace-Text "32×84=" "44×93="
Replace-Text "59×73=" "16×49="
Replace-Text "30×61=" "45×77="
Replace-Text "20×82=" "56×30="

Replace-Text "88×74=" "68×96="
Replace-Text "90×54=" "63×63="
Replace-Text "50×21=" "14×53="
Replace-Text "54×52=" "71×59="
Replace-Text "59×25=" "41×18="

Replace-Text "79×76=" "29×62="
Replace-Text "97×14=" "87×83="
Replace-Text "88×86=" "90×57="
Replace-Text "96×57=" "27×59="
Replace-Text "91×46=" "52×14="
